$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (first data row)
$t.Cell(1, 1).Range.Text = "39÷2=19, 1"
$t.Cell(1, 3).Range.Text = "28÷5=5, 3"
$t.Cell(1, 4).Range.Text = "62÷3=20, 2"
$t.Cell(1, 5).Range.Text = "86÷6=14, 2"

# Row 5 (second data row)
$t.Cell(5, 1).Range.Text = "73÷7=10, 3"
$t.Cell(5, 2).Range.Text = "18÷6=3, 0"
$t.Cell(5, 3).Range.Text = "47÷4=11, 3"
$t.Cell(5, 4).Range.Text = "44÷2=22, 0"
$t.Cell(5, 5).Range.Text = "25÷4=6, 1"

# Row 9 (third data row)
$t.Cell(9, 1).Range.Text = "41÷4=10, 1"
$t.Cell(9, 2).Range.Text = "77÷6=12, 5"
$t.Cell(9, 3).Range.Text = "14÷7=2, 0"
$t.Cell(9, 4).Range.Text = "68÷6=11, 2"
$t.Cell(9, 5).Range.Text = "30÷5=6, 0"

# Row 13 (fourth data row)
$t.Cell(13, 1).Range.Text = "20÷7=2, 6"
$t.Cell(13, 2).Range.Text = "33÷6=5, 3"
$t.Cell(13, 3).Range.Text = "60÷4=15, 0"
$t.Cell(13, 4).Range.Text = "22÷2=11, 0"
$t.Cell(13, 5).Range.Text = "79÷3=26, 1"

# Row 17 (fifth data row)
$t.Cell(17, 1).Range.Text = "98÷9=10, 8"
$t.Cell(17, 2).Range.Text = "53÷7=7, 4"
$t.Cell(17, 3).Range.Text = "55÷2=27, 1"
$t.Cell(17, 4).Range.Text = "86÷8=10, 6"
$t.Cell(17, 5).Range.Text = "21÷3=7, 0"
